$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header for the electrode location column, matching the
# existing header formatting (bold, centered, bordered)
$ws.Cells.Item(1, 3).Value = "Electrode Locations"
$ws.Range("A1").Copy()
$ws.Range("C1").PasteSpecial(-4122)

# Sorted (A1-O15) electrode data: filename suffix, unnormalized Pmax, electrode location
$data = @(
    @{Loc="A2"; Val=58.856285},
    @{Loc="A3"; Val=54.625608},
    @{Loc="A6"; Val=35.001385},
    @{Loc="A8"; Val=45.29409},
    @{Loc="A9"; Val=49.754141},
    @{Loc="A11"; Val=53.216597},
    @{Loc="B14"; Val=50.576975},
    @{Loc="B15"; Val=51.494471},
    @{Loc="C1"; Val=58.128114},
    @{Loc="C3"; Val=53.660781},
    @{Loc="C9"; Val=49.916159},
    @{Loc="C13"; Val=48.255928},
    @{Loc="C15"; Val=47.609676},
    @{Loc="D6"; Val=44.35839},
    @{Loc="D8"; Val=44.917261},
    @{Loc="D12"; Val=47.172773},
    @{Loc="E1"; Val=61.72164},
    @{Loc="E4"; Val=50.405854},
    @{Loc="E6"; Val=48.235903},
    @{Loc="E7"; Val=50.744454},
    @{Loc="E10"; Val=49.553894},
    @{Loc="E11"; Val=51.281481},
    @{Loc="E15"; Val=48.907641},
    @{Loc="F14"; Val=52.42471},
    @{Loc="G1"; Val=63.827876},
    @{Loc="G4"; Val=53.560658},
    @{Loc="G5"; Val=49.974413},
    @{Loc="G7"; Val=52.344611},
    @{Loc="G9"; Val=51.547263},
    @{Loc="G12"; Val=49.967131},
    @{Loc="G13"; Val=48.789314},
    @{Loc="G15"; Val=51.350657},
    @{Loc="I3"; Val=53.500583},
    @{Loc="I5"; Val=54.793088},
    @{Loc="I8"; Val=53.333104},
    @{Loc="I10"; Val=50.38765},
    @{Loc="I11"; Val=47.107237},
    @{Loc="I15"; Val=51.44714},
    @{Loc="J2"; Val=57.101392},
    @{Loc="J13"; Val=50.058152},
    @{Loc="K1"; Val=58.210033},
    @{Loc="K3"; Val=58.122652},
    @{Loc="K5"; Val=56.983064},
    @{Loc="K10"; Val=48.871233},
    @{Loc="K12"; Val=48.620014},
    @{Loc="K14"; Val=51.714743},
    @{Loc="K15"; Val=52.386481},
    @{Loc="L6"; Val=48.883976},
    @{Loc="M2"; Val=60.177916},
    @{Loc="M5"; Val=57.865972},
    @{Loc="M7"; Val=50.462288},
    @{Loc="M10"; Val=44.751602},
    @{Loc="M11"; Val=51.674693},
    @{Loc="M13"; Val=50.43134},
    @{Loc="M15"; Val=51.456242},
    @{Loc="N3"; Val=59.147554},
    @{Loc="O2"; Val=57.510988},
    @{Loc="O3"; Val=58.83262},
    @{Loc="O5"; Val=57.430889},
    @{Loc="O8"; Val=49.998078},
    @{Loc="O10"; Val=46.684898},
    @{Loc="O14"; Val=43.790416}
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = "$($item.Loc)_bipolar_10V_100kHz.txt"
    $ws.Cells.Item($row, 2).Value = $item.Val
    $ws.Cells.Item($row, 3).Value = $item.Loc
    $row++
}

